$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 42.32036466666667
$ws.Range("H2").Value = 126.961094
$ws.Range("I2").Value = 0.285778576657872
$ws.Range("J2").Value = 0.2880046678857171
$ws.Range("M2").Value = 1.378421333333333
$ws.Range("N2").Value = 4.135264
$ws.Range("O2").Value = 0.01656231489052403
$ws.Range("P2").Value = 0.01794267551419991
$ws.Range("Q2").Value = 58.33529349097957
$ws.Range("R2").Value = 525.0176414188161
$ws.Range("S2").Value = 0.004733154775573438
$ws.Range("T2").Value = 0.005167574302448334
$ws.Range("G3").Value = 42.32036466666667
$ws.Range("H3").Value = 126.961094
$ws.Range("I3").Value = 0.285778576657872
$ws.Range("J3").Value = 0.2880046678857171
$ws.Range("O3").Value = 0.2170932623988173
$ws.Range("P3").Value = 0.2351865659654651
$ws.Range("Q3").Value = 764.6394396350322
$ws.Range("R3").Value = 6881.754956715289
$ws.Range("S3").Value = 0.06204060353034793
$ws.Range("T3").Value = 0.06773482882206608
$ws.Range("G4").Value = 42.32036466666667
$ws.Range("H4").Value = 126.961094
$ws.Range("I4").Value = 0.285778576657872
$ws.Range("J4").Value = 0.2880046678857171
$ws.Range("M4").Value = 17.58286933333333
$ws.Range("N4").Value = 52.748608
$ws.Range("O4").Value = 0.2112656061941426
$ws.Range("P4").Value = 0.22887321273073
$ws.Range("Q4").Value = 744.1134420730168
$ws.Range("R4").Value = 6697.020978657152
$ws.Range("S4").Value = 0.06037518423492459
$ws.Range("T4").Value = 0.06591655362045098
$ws.Range("G5").Value = 42.32036466666667
$ws.Range("H5").Value = 126.961094
$ws.Range("I5").Value = 0.285778576657872
$ws.Range("J5").Value = 0.2880046678857171
$ws.Range("M5").Value = 19.2082395
$ws.Range("N5").Value = 38.416479
$ws.Range("O5").Value = 0.2307951156866419
$ws.Range("P5").Value = 0.1666869194070983
$ws.Range("Q5").Value = 812.8997002446711
$ws.Range("R5").Value = 4877.398201468027
$ws.Range("S5").Value = 0.06595629966051744
$ws.Range("T5").Value = 0.04800661086473465
$ws.Range("G6").Value = 42.32036466666667
$ws.Range("H6").Value = 126.961094
$ws.Range("I6").Value = 0.285778576657872
$ws.Range("J6").Value = 0.2880046678857171
$ws.Range("M6").Value = 26.988955
$ws.Range("N6").Value = 80.966865
$ws.Range("O6").Value = 0.3242837008298742
$ws.Range("P6").Value = 0.3513106263825066
$ws.Range("Q6").Value = 1142.182417572257
$ws.Range("R6").Value = 10279.64175815031
$ws.Range("S6").Value = 0.09267333445650865
$ws.Range("T6").Value = 0.1011791002760171
$ws.Range("I7").Value = 0.04213668412459876
$ws.Range("J7").Value = 0.04246491062777905
$ws.Range("M7").Value = 1.378421333333333
$ws.Range("N7").Value = 4.135264
$ws.Range("O7").Value = 0.01656231489052403
$ws.Range("P7").Value = 0.01794267551419991
$ws.Range("Q7").Value = 8.601259982087111
$ws.Range("R7").Value = 77.411339838784
$ws.Range("S7").Value = 0.0006978810309141496
$ws.Range("T7").Value = 0.0007619341121337387
$ws.Range("I8").Value = 0.04213668412459876
$ws.Range("J8").Value = 0.04246491062777905
$ws.Range("O8").Value = 0.2170932623988173
$ws.Range("P8").Value = 0.2351865659654651
$ws.Range("S8").Value = 0.009147590223277597
$ws.Range("T8").Value = 0.009987176504577738
$ws.Range("I9").Value = 0.04213668412459876
$ws.Range("J9").Value = 0.04246491062777905
$ws.Range("M9").Value = 17.58286933333333
$ws.Range("N9").Value = 52.748608
$ws.Range("O9").Value = 0.2112656061941426
$ws.Range("P9").Value = 0.22887321273073
$ws.Range("Q9").Value = 109.7159676144498
$ws.Range("R9").Value = 987.4437085300478
$ws.Range("S9").Value = 0.008902032114594461
$ws.Range("T9").Value = 0.009719080523703112
$ws.Range("I10").Value = 0.04213668412459876
$ws.Range("J10").Value = 0.04246491062777905
$ws.Range("M10").Value = 19.2082395
$ws.Range("N10").Value = 38.416479
$ws.Range("O10").Value = 0.2307951156866419
$ws.Range("P10").Value = 0.1666869194070983
$ws.Range("Q10").Value = 119.858172347179
$ws.Range("R10").Value = 719.149034083074
$ws.Range("S10").Value = 0.009724940887188258
$ws.Range("T10").Value = 0.00707834513544224
$ws.Range("I11").Value = 0.04213668412459876
$ws.Range("J11").Value = 0.04246491062777905
$ws.Range("M11").Value = 26.988955
$ws.Range("N11").Value = 80.966865
$ws.Range("O11").Value = 0.3242837008298742
$ws.Range("P11").Value = 0.3513106263825066
$ws.Range("Q11").Value = 168.4093339142433
$ws.Range("R11").Value = 1515.68400522819
$ws.Range("S11").Value = 0.01366423986862429
$ws.Range("T11").Value = 0.01491837435192222
$ws.Range("G12").Value = 42.241047
$ws.Range("H12").Value = 126.723141
$ws.Range("I12").Value = 0.2852429647825406
$ws.Range("J12").Value = 0.2874648838260633
$ws.Range("M12").Value = 1.378421333333333
$ws.Range("N12").Value = 4.135264
$ws.Range("O12").Value = 0.01656231489052403
$ws.Range("P12").Value = 0.01794267551419991
$ws.Range("Q12").Value = 58.22596032713601
$ws.Range("R12").Value = 524.0336429442241
$ws.Range("S12").Value = 0.004724283803035093
$ws.Range("T12").Value = 0.005157889132218228
$ws.Range("G13").Value = 42.241047
$ws.Range("H13").Value = 126.723141
$ws.Range("I13").Value = 0.2852429647825406
$ws.Range("J13").Value = 0.2874648838260633
$ws.Range("O13").Value = 0.2170932623988173
$ws.Range("P13").Value = 0.2351865659654651
$ws.Range("Q13").Value = 763.2063372345482
$ws.Range("R13").Value = 6868.857035110933
$ws.Range("S13").Value = 0.06192432580095267
$ws.Range("T13").Value = 0.0676078788627132
$ws.Range("G14").Value = 42.241047
$ws.Range("H14").Value = 126.723141
$ws.Range("I14").Value = 0.2852429647825406
$ws.Range("J14").Value = 0.2874648838260633
$ws.Range("M14").Value = 17.58286933333333
$ws.Range("N14").Value = 52.748608
$ws.Range("O14").Value = 0.2112656061941426
$ws.Range("P14").Value = 0.22887321273073
$ws.Range("Q14").Value = 742.718809904192
$ws.Range("R14").Value = 6684.469289137727
$ws.Range("S14").Value = 0.0602620278673979
$ws.Range("T14").Value = 0.06579301150853717
$ws.Range("G15").Value = 42.241047
$ws.Range("H15").Value = 126.723141
$ws.Range("I15").Value = 0.2852429647825406
$ws.Range("J15").Value = 0.2874648838260633
$ws.Range("M15").Value = 19.2082395
$ws.Range("N15").Value = 38.416479
$ws.Range("O15").Value = 0.2307951156866419
$ws.Range("P15").Value = 0.1666869194070983
$ws.Range("Q15").Value = 811.3761475067566
$ws.Range("R15").Value = 4868.25688504054
$ws.Range("S15").Value = 0.06583268305578717
$ws.Range("T15").Value = 0.0479166359226859
$ws.Range("G16").Value = 42.241047
$ws.Range("H16").Value = 126.723141
$ws.Range("I16").Value = 0.2852429647825406
$ws.Range("J16").Value = 0.2874648838260633
$ws.Range("M16").Value = 26.988955
$ws.Range("N16").Value = 80.966865
$ws.Range("O16").Value = 0.3242837008298742
$ws.Range("P16").Value = 0.3513106263825066
$ws.Range("Q16").Value = 1140.041716635885
$ws.Range("R16").Value = 10260.37544972296
$ws.Range("S16").Value = 0.09249964425536772
$ws.Range("T16").Value = 0.1009894683999088
$ws.Range("G17").Value = 3.4338745
$ws.Range("H17").Value = 6.867749
$ws.Range("I17").Value = 0.02318807445921414
$ws.Range("J17").Value = 0.0155791330048516
$ws.Range("M17").Value = 1.378421333333333
$ws.Range("N17").Value = 4.135264
$ws.Range("O17").Value = 0.01656231489052403
$ws.Range("P17").Value = 0.01794267551419991
$ws.Range("Q17").Value = 4.733325866789333
$ws.Range("R17").Value = 28.399955200736
$ws.Range("S17").Value = 0.0003840481908984223
$ws.Range("T17").Value = 0.0002795313282986144
$ws.Range("G18").Value = 3.4338745
$ws.Range("H18").Value = 6.867749
$ws.Range("I18").Value = 0.02318807445921414
$ws.Range("J18").Value = 0.0155791330048516
$ws.Range("O18").Value = 0.2170932623988173
$ws.Range("P18").Value = 0.2351865659654651
$ws.Range("Q18").Value = 62.04284613655801
$ws.Range("R18").Value = 372.257076819348
$ws.Range("S18").Value = 0.005033974733097487
$ws.Range("T18").Value = 0.003664002792130285
$ws.Range("G19").Value = 3.4338745
$ws.Range("H19").Value = 6.867749
$ws.Range("I19").Value = 0.02318807445921414
$ws.Range("J19").Value = 0.0155791330048516
$ws.Range("M19").Value = 17.58286933333333
$ws.Range("N19").Value = 52.748608
$ws.Range("O19").Value = 0.2112656061941426
$ws.Range("P19").Value = 0.22887321273073
$ws.Range("Q19").Value = 60.37736664056533
$ws.Range("R19").Value = 362.264199843392
$ws.Range("S19").Value = 0.00489884260710079
$ws.Range("T19").Value = 0.003565646222379737
$ws.Range("G20").Value = 3.4338745
$ws.Range("H20").Value = 6.867749
$ws.Range("I20").Value = 0.02318807445921414
$ws.Range("J20").Value = 0.0155791330048516
$ws.Range("M20").Value = 19.2082395
$ws.Range("N20").Value = 38.416479
$ws.Range("O20").Value = 0.2307951156866419
$ws.Range("P20").Value = 0.1666869194070983
$ws.Range("Q20").Value = 65.95868380894275
$ws.Range("R20").Value = 263.834735235771
$ws.Range("S20").Value = 0.005351694327364794
$ws.Range("T20").Value = 0.002596837687612164
$ws.Range("G21").Value = 3.4338745
$ws.Range("H21").Value = 6.867749
$ws.Range("I21").Value = 0.02318807445921414
$ws.Range("J21").Value = 0.0155791330048516
$ws.Range("M21").Value = 26.988955
$ws.Range("N21").Value = 80.966865
$ws.Range("O21").Value = 0.3242837008298742
$ws.Range("P21").Value = 0.3513106263825066
$ws.Range("Q21").Value = 92.6766843561475
$ws.Range("R21").Value = 556.060106136885
$ws.Range("S21").Value = 0.007519514600752644
$ws.Range("T21").Value = 0.005473114974430796
$ws.Range("G22").Value = 53.85273233333334
$ws.Range("H22").Value = 161.558197
$ws.Range("I22").Value = 0.3636536999757743
$ws.Range("J22").Value = 0.3664864046555889
$ws.Range("M22").Value = 1.378421333333333
$ws.Range("N22").Value = 4.135264
$ws.Range("O22").Value = 0.01656231489052403
$ws.Range("P22").Value = 0.01794267551419991
$ws.Range("Q22").Value = 74.23175510655646
$ws.Range("R22").Value = 668.0857959590081
$ws.Range("S22").Value = 0.006022947090102926
$ws.Range("T22").Value = 0.006575746639100996
$ws.Range("G23").Value = 53.85273233333334
$ws.Range("H23").Value = 161.558197
$ws.Range("I23").Value = 0.3636536999757743
$ws.Range("J23").Value = 0.3664864046555889
$ws.Range("O23").Value = 0.2170932623988173
$ws.Range("P23").Value = 0.2351865659654651
$ws.Range("Q23").Value = 973.0049208817162
$ws.Range("R23").Value = 8757.044287935445
$ws.Range("S23").Value = 0.07894676811114155
$ws.Range("T23").Value = 0.08619267898397781
$ws.Range("G24").Value = 53.85273233333334
$ws.Range("H24").Value = 161.558197
$ws.Range("I24").Value = 0.3636536999757743
$ws.Range("J24").Value = 0.3664864046555889
$ws.Range("M24").Value = 17.58286933333333
$ws.Range("N24").Value = 52.748608
$ws.Range("O24").Value = 0.2112656061941426
$ws.Range("P24").Value = 0.22887321273073
$ws.Range("Q24").Value = 946.8855558599751
$ws.Range("R24").Value = 8521.970002739776
$ws.Range("S24").Value = 0.07682751937012482
$ws.Range("T24").Value = 0.08387892085565901
$ws.Range("G25").Value = 53.85273233333334
$ws.Range("H25").Value = 161.558197
$ws.Range("I25").Value = 0.3636536999757743
$ws.Range("J25").Value = 0.3664864046555889
$ws.Range("M25").Value = 19.2082395
$ws.Range("N25").Value = 38.416479
$ws.Range("O25").Value = 0.2307951156866419
$ws.Range("P25").Value = 0.1666869194070983
$ws.Range("Q25").Value = 1034.416180388061
$ws.Range("R25").Value = 6206.497082328364
$ws.Range("S25").Value = 0.0839294977557842
$ws.Range("T25").Value = 0.06108848979662338
$ws.Range("G26").Value = 53.85273233333334
$ws.Range("H26").Value = 161.558197
$ws.Range("I26").Value = 0.3636536999757743
$ws.Range("J26").Value = 0.3664864046555889
$ws.Range("M26").Value = 26.988955
$ws.Range("N26").Value = 80.966865
$ws.Range("O26").Value = 0.3242837008298742
$ws.Range("P26").Value = 0.3513106263825066
$ws.Range("Q26").Value = 1453.428969571379
$ws.Range("R26").Value = 13080.86072614241
$ws.Range("S26").Value = 0.1179269676486208
$ws.Range("T26").Value = 0.1287505683802277
